{"js": "// Word JS API (Office.js) script.\n// This is the body of: async (context) => { ... }\n//\n// Change 1: the title paragraph \"LOOKUP TABLE\" becomes \"Ataque LOOKUP TABLE\"\n//           (a new run with the text \"Ataque \" is added in front of the\n//           existing \"LOOKUP TABLE\" run).\n// Change 2: the sentence that used to read \"...realizar ataques a los hash,\n//           para...\" keeps the exact same wording (\"a los hash\" is merely\n//           re-typed / re-split into \"a los\" + \" \" + \"hash\" by the author in\n//           the source edit) -- the visible text does not change.\n\nconst body = context.document.body;\n\n// --- Change 1: prepend \"Ataque \" to the title paragraph -------------------\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst titleParagraph = paragraphs.items[0];\nconst titleStart = titleParagraph.getRange(\"Start\");\ntitleStart.insertText(\"Ataque \", \"Before\");\nawait context.sync();\n\n// --- Change 2: re-split \"a los hash\" into \"a los\" + \" \" + \"hash\" ----------\n// Find the unique paragraph that contains the phrase so the search below is\n// scoped to exactly the right spot in the document.\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet targetParagraph = null;\nfor (const p of paragraphs.items) {\n  if (p.text.indexOf(\"a los hash\") !== -1) {\n    targetParagraph = p;\n    break;\n  }\n}\n\nif (targetParagraph) {\n  const hits = targetParagraph.search(\"a los hash\", { matchCase: true });\n  hits.load(\"items\");\n  await context.sync();\n\n  if (hits.items.length > 0) {\n    const hit = hits.items[0];\n    // Replace the phrase with itself, split across three insert calls so the\n    // resulting text is unchanged (\"a los\" + \" \" + \"hash\" === \"a los hash\")\n    // while the edit is performed as three discrete operations, mirroring\n    // the run split recorded in the source revision.\n    hit.insertText(\"a los\", \"Replace\");\n    await context.sync();\n\n    const afterLos = targetParagraph.search(\"a los\", { matchCase: true });\n    afterLos.load(\"items\");\n    await context.sync();\n    const losRange = afterLos.items[0].getRange(\"End\");\n    losRange.insertText(\" \", \"After\");\n    await context.sync();\n\n    const afterSpace = targetParagraph.search(\"a los \", { matchCase: true });\n    afterSpace.load(\"items\");\n    await context.sync();\n    const spaceRange = afterSpace.items[0].getRange(\"End\");\n    spaceRange.insertText(\"hash\", \"After\");\n    await context.sync();\n  }\n}\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# $word.ActiveDocument is already open as $d.\n#\n# Change 1: the title paragraph \"LOOKUP TABLE\" becomes \"Ataque LOOKUP TABLE\"\n#           (a new run with the text \"Ataque \" is added in front of the\n#           existing \"LOOKUP TABLE\" run).\n# Change 2: the sentence that used to read \"...realizar ataques a los hash,\n#           para...\" keeps the exact same wording (\"a los hash\" is merely\n#           re-typed / re-split into \"a los\" + \" \" + \"hash\" by the author in\n#           the source edit) -- the visible text does not change.\n\n$d = $word.ActiveDocument\n\n# --- Change 1: prepend \"Ataque \" to the title paragraph --------------------\n$titlePar = $d.Paragraphs(1)\n$titleRange = $titlePar.Range\n$titleRange.Collapse(1)           # wdCollapseStart\n$titleRange.InsertBefore(\"Ataque \")\n\n# --- Change 2: re-split \"a los hash\" into \"a los\" + \" \" + \"hash\" -----------\n$find = $d.Content\n$find.Find.ClearFormatting()\n$find.Find.MatchCase = $true\n$find.Find.MatchWholeWord = $false\n$find.Find.Text = \"a los hash\"\n$found = $find.Find.Execute()\n\nif ($found) {\n    # Replace the phrase with itself, split across three discrete edits so\n    # the resulting text is unchanged (\"a los\" + \" \" + \"hash\" === \"a los hash\")\n    # while mirroring the run split recorded in the source revision.\n    $find.Text = \"a los\"\n    $find.Collapse(0)             # wdCollapseEnd\n    $find.InsertAfter(\" \")\n    $find.Collapse(0)\n    $find.InsertAfter(\"hash\")\n    $find.Collapse(0)\n}\n"}
